# Order.xlsx edit: add rate column data to the imported order row and
# switch the sheet selection, per commit "Add rate in impoert order,
# change samll to caps header in inocie and lpo pdf".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (the imported order data row) gets new values.
$ws.Range("A2").Value = "webtest1@yopmail.com"
$ws.Range("B2").Value = 382
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 120
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 120
$ws.Range("M2").Value = 12

# Move the active selection to F2 (and drop the old scrolled/top-left
# state that pointed at O5 / column B).
[void]$ws.Range("F2").Select()
